# Applies crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.318.74"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.869.48"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'236.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D7").Value = "'0.4710"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("E8").Value = "  +2.16%  "
$ws.Range("D9").Value = "'0.06623"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("D10").Value = "'21.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "'97.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "1.876.47"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "'5.157"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "'0.6892"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").Value = "'274.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "30.315.17"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "'14.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.96%  "
$ws.Range("D19").Value = "'0.000007725"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.70%  "
$ws.Range("D20").Value = "'1.0000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").Value = "2.122.93"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").Value = "'5.317"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'6.225"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").Value = "'167.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("D26").Value = "'9.282"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("D27").Value = "'18.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("D29").Value = "'1.369"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("E30").Value = "  +2.32%  "
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("D33").Value = "'4.093"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("D34").Value = "'0.04705"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("D35").Value = "'1.134"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Value = "'0.7037"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Value = "'2.703"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").Value = "'0.01882"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("D39").Value = "'2.637"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.48%  "
$ws.Range("D40").Value = "'6.311"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").Value = "'73.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.65%  "
$ws.Range("D42").Value = "'1.965"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.8425"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.4172"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").Value = "'0.9994"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").Value = "'103.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("D47").Value = "'7.141"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").Value = "'9.274"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("D49").Value = "'932.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.01%  "
$ws.Range("D50").Value = "'34.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("D51").Value = "'0.05666"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.42%  "
